# 1st testcase of forgotPassword scenario added
#
# The author marked the Registration_Account_Setup6 test case ("user4@gmail.com")
# as passing: the Runmode/result cell on that scenario's sheet flips from
# "fail" to "pass", and the corresponding roll-up row on the "Test Cases"
# summary sheet flips from "Fail" to "Pass". Selection/active-sheet focus
# moves from the Registration_Account_Setup6 sheet (B7) to the Test Cases
# sheet (E8), which is where that roll-up result lives.

$wb = $excel.ActiveWorkbook

# --- Registration_Account_Setup6: flip the scenario's result to "pass" ---
$wsSetup6 = $wb.Worksheets.Item("Registration_Account_Setup6")
$wsSetup6.Range("F2").Value = "pass"

# --- Test Cases: flip the matching roll-up Results cell to "Pass" ---
$wsTests = $wb.Worksheets.Item("Test Cases")
$wsTests.Range("E8").Value = "Pass"

# --- Move focus to the Test Cases sheet / the cell that was just updated ---
$wsTests.Activate()
$wsTests.Range("E8").Select()
